$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the data of row 3 and row 4 (the two species records traded
# places). Only the columns whose values actually differ between the two
# rows need to be touched: A, B, D, E, F, G, H, M, Q, R, S.
# (Column M - "Aktivitet" - only had a value on row 4 and moves to row 3.)

$cols = @("A", "B", "D", "E", "F", "G", "H", "M", "Q", "R", "S")

$row3 = @{}
$row4 = @{}

foreach ($col in $cols) {
    $row3[$col] = $ws.Range("${col}3").Value2
    $row4[$col] = $ws.Range("${col}4").Value2
}

foreach ($col in $cols) {
    $ws.Range("${col}3").Value2 = $row4[$col]
    $ws.Range("${col}4").Value2 = $row3[$col]
}
